$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1457.7028
$ws.Range("I129").Value = 819.6
$ws.Range("J129").Value = 1557.4062
$ws.Range("K129").Value = 2458.8
$ws.Range("L129").Value = 4672.2186
$ws.Range("M129").Value = 2541.2
$ws.Range("N129").Value = -14672.2186

$ws.Range("H137").Value = 8334955.5
$ws.Range("I137").Value = 13515231
$ws.Range("K137").Value = 40545693
$ws.Range("M137").Value = -40543143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4515592
$ws.Range("I32").Value = 6250.254
$ws.Range("J32").Value = 30341822
$ws.Range("K32").Value = 6250.254
$ws.Range("L32").Value = 30341822
$ws.Range("M32").Value = -5963.254
$ws.Range("N32").Value = -30342396

$ws.Range("H61").Value = 2263.6155
$ws.Range("I61").Value = 1518.8125
$ws.Range("J61").Value = 2781.739
$ws.Range("K61").Value = 1518.8125
$ws.Range("L61").Value = 2781.739
$ws.Range("M61").Value = -1306.8125
$ws.Range("N61").Value = -3205.739

$ws.Range("H110").Value = 2434.24
$ws.Range("I110").Value = 1516.3334
$ws.Range("K110").Value = 1516.3334
$ws.Range("M110").Value = 528.6666

$ws.Range("H136").Value = 2263.6155
$ws.Range("I136").Value = 1518.8125
$ws.Range("J136").Value = 2781.739
$ws.Range("K136").Value = 4556.4375
$ws.Range("L136").Value = 8345.217000000001
$ws.Range("M136").Value = -2006.4375
$ws.Range("N136").Value = -13445.217

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4085.5088
$ws.Range("I134").Value = 1703.4333
$ws.Range("J134").Value = 6732.2593
$ws.Range("K134").Value = 5110.2999
$ws.Range("L134").Value = 20196.7779
$ws.Range("M134").Value = -2575.2999
$ws.Range("N134").Value = -25266.7779

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2428.111
$ws.Range("I132").Value = 1637.5
$ws.Range("J132").Value = 3060.6
$ws.Range("K132").Value = 4912.5
$ws.Range("L132").Value = 9181.799999999999
$ws.Range("M132").Value = -2382.5
$ws.Range("N132").Value = -14241.8

$ws.Range("H134").Value = 2033.5869
$ws.Range("I134").Value = 992.2414
$ws.Range("J134").Value = 3810
$ws.Range("K134").Value = 2976.7242
$ws.Range("L134").Value = 11430
$ws.Range("M134").Value = -441.7242000000001
$ws.Range("N134").Value = -16500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 509.16666
$ws.Range("I2").Value = 53.8
$ws.Range("J2").Value = 834.4286
$ws.Range("K2").Value = 322.8
$ws.Range("L2").Value = 5006.571599999999
$ws.Range("M2").Value = -209.8
$ws.Range("N2").Value = -5232.571599999999

$ws.Range("H15").Value = 320.9
$ws.Range("I15").Value = 320.9
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 962.6999999999999
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -822.6999999999999
$ws.Range("N15").ClearContents()

$ws.Range("H17").Value = 409.83334
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 431.8
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 1295.4
$ws.Range("M17").Value = -731
$ws.Range("N17").Value = -1633.4

$ws.Range("H22").Value = 3625.5
$ws.Range("I22").Value = 3834
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 11502
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = -11333
$ws.Range("N22").Value = -9338

$ws.Range("H24").Value = 4500
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 4500
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 13500
$ws.Range("N24").Value = -13960
$ws.Range("M24").ClearContents()

$ws.Range("H27").Value = 3625.5
$ws.Range("I27").Value = 3834
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 11502
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = -11400
$ws.Range("N27").Value = -9204

$ws.Range("H32").Value = 2500
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 3500
$ws.Range("K32").Value = 4500
$ws.Range("L32").Value = 10500
$ws.Range("M32").Value = -4217
$ws.Range("N32").Value = -11066

$ws.Range("H39").Value = 10000
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -30588

$ws.Range("H58").Value = 2475
$ws.Range("I58").Value = 2475
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 7425
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -7297
$ws.Range("N58").ClearContents()

$ws.Range("H76").Value = 3353.7778
$ws.Range("I76").Value = 1930.6666
$ws.Range("J76").Value = 6200
$ws.Range("K76").Value = 5791.9998
$ws.Range("L76").Value = 18600
$ws.Range("M76").Value = -5408.9998
$ws.Range("N76").Value = -19366

$ws.Range("H79").Value = 3353.7778
$ws.Range("I79").Value = 1930.6666
$ws.Range("J79").Value = 6200
$ws.Range("K79").Value = 5791.9998
$ws.Range("L79").Value = 18600
$ws.Range("M79").Value = -4465.9998
$ws.Range("N79").Value = -21252

$ws.Range("H93").Value = 5309.4
$ws.Range("J93").Value = 5309.4
$ws.Range("L93").Value = 15928.2
$ws.Range("N93").Value = -19672.2

$ws.Range("H100").Value = 3824.75
$ws.Range("J100").Value = 3824.75
$ws.Range("L100").Value = 11474.25
$ws.Range("N100").Value = -13096.25

$ws.Range("H106").Value = 3863.3635
$ws.Range("J106").Value = 3863.3635
$ws.Range("L106").Value = 11590.0905
$ws.Range("N106").Value = -13482.0905

$ws.Range("H112").Value = 32525004
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 32525004
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 97575012
$ws.Range("N112").Value = -97577228
$ws.Range("M112").ClearContents()

$ws.Range("H123").Value = 3733.1667
$ws.Range("I123").Value = 1999.5
$ws.Range("J123").Value = 4600
$ws.Range("K123").Value = 5998.5
$ws.Range("L123").Value = 13800
$ws.Range("M123").Value = -3548.5
$ws.Range("N123").Value = -18700

$ws.Range("H125").Value = 6203.7036
$ws.Range("I125").Value = 2650
$ws.Range("J125").Value = 6488
$ws.Range("K125").Value = 7950
$ws.Range("L125").Value = 19464
$ws.Range("M125").Value = -3030
$ws.Range("N125").Value = -29304

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2930.1428
$ws.Range("I113").Value = 2702.2
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 2702.2
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -532.1999999999998
$ws.Range("N113").Value = -7840

$ws.Range("H126").Value = 31252648
$ws.Range("I126").Value = 50001640
$ws.Range("J126").Value = 4326.6665
$ws.Range("K126").Value = 150004920
$ws.Range("L126").Value = 12979.9995
$ws.Range("M126").Value = -150002450
$ws.Range("N126").Value = -17919.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3130.625
$ws.Range("I68").Value = 2740
$ws.Range("K68").Value = 2740
$ws.Range("M68").Value = -1991

$ws.Range("H71").Value = 3130.625
$ws.Range("I71").Value = 2740
$ws.Range("K71").Value = 13700
$ws.Range("M71").Value = -9956

$ws.Range("H93").Value = 228801.14
$ws.Range("I93").Value = 386051.62
$ws.Range("J93").Value = 1661.5555
$ws.Range("K93").Value = 386051.62
$ws.Range("L93").Value = 1661.5555
$ws.Range("M93").Value = -384803.62
$ws.Range("N93").Value = -4157.5555

$ws.Range("H132").Value = 28604954
$ws.Range("I132").Value = 34521364
$ws.Range("J132").Value = 8966.333000000001
$ws.Range("K132").Value = 103564092
$ws.Range("L132").Value = 26898.999
$ws.Range("M132").Value = -103561562
$ws.Range("N132").Value = -31958.999

$ws.Range("H136").Value = 17858270
$ws.Range("I136").Value = 26316626
$ws.Range("J136").Value = 1741.5555
$ws.Range("K136").Value = 78949878
$ws.Range("L136").Value = 5224.666499999999
$ws.Range("M136").Value = -78947328
$ws.Range("N136").Value = -10324.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8671.182000000001
$ws.Range("I122").Value = 9912.875
$ws.Range("J122").Value = 5360
$ws.Range("K122").Value = 29738.625
$ws.Range("L122").Value = 16080
$ws.Range("M122").Value = -27288.625
$ws.Range("N122").Value = -20980

$ws.Range("H132").Value = 2667.4736
$ws.Range("I132").Value = 2208.25
$ws.Range("J132").Value = 3001.4546
$ws.Range("K132").Value = 6624.75
$ws.Range("L132").Value = 9004.363799999999
$ws.Range("M132").Value = -4094.75
$ws.Range("N132").Value = -14064.3638
